$d = $word.ActiveDocument

# --- Change 1: "16-99" -> "18" + "+" (two runs) and move the _GoBack bookmark here ---

# Remove the existing (old) _GoBack bookmark first, since it will be re-created
# at its new location further below (bookmark names must stay unique).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$rng = $d.Content
$rng.Find.Execute("16-99", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Text = "18"
$rng.Collapse(0)
$rng.InsertAfter("+")
$rng.Font.Name = "Times New Roman"
$rng.Font.NameBi = "Times New Roman"

# Place a fresh collapsed bookmark right after "+" (before the paragraph mark).
# Placing a bookmark collapsed exactly at the paragraph-mark boundary is flaky,
# so insert a throwaway character, bookmark across it, then remove the
# character again - this leaves a clean, correctly collapsed bookmark behind.
$rng.Collapse(0)
$rng.InsertAfter("X")
$d.Bookmarks.Add("_GoBack", $rng)
$rng.Text = ""

# --- Change 2: fix typo "безуспкешно" -> "безуспешно" ---
# Grab the whole run that contains the typo (not just the matched word) and
# rewrite its text in one shot. A plain text replace on an identically
# formatted neighbour run causes this engine to auto-merge adjacent runs
# that share formatting, which would collapse several runs into one and
# break the run structure the diff expects. Briefly toggling a formatting
# property (Bold) around the text edit suppresses that auto-merge, and then
# restoring the original value leaves the run's formatting untouched.
$rng2 = $d.Content
$rng2.Find.Execute("ийся безуспкешно предотвратить ")
$origBold2 = $rng2.Font.Bold
$rng2.Font.Bold = 1 - $origBold2
$rng2.Text = "ийся безуспешно предотвратить "
$rng2.Font.Bold = $origBold2
